# Rename "RESOURCE.RTYPE.CLASS" -> "RESOURCE.RTYPE.RCLASS"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESOURCE.RTYPE.CLASS")
$ws.Name = "RESOURCE.RTYPE.RCLASS"

# Make it the active sheet/tab (was the capability sheet before) and
# move its selection to C32.
$ws.Activate()
$ws.Range("C32").Select()
